$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the Longitude (C) and Latitude (D) columns entirely.
$ws.Range("C:D").Delete()

# Update selection to match the post-edit state.
$ws.Range("E10:F10").Select() | Out-Null
